$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(2,1).Value = 0
$ws.Cells.Item(2,2).Value = 73
$ws.Cells.Item(3,1).Value = 1
$ws.Cells.Item(3,2).Value = 76
$ws.Cells.Item(4,1).Value = 2
$ws.Cells.Item(4,2).Value = 79
$ws.Cells.Item(5,1).Value = 3
$ws.Cells.Item(5,2).Value = 82
$ws.Cells.Item(6,1).Value = 4
$ws.Cells.Item(6,2).Value = 85
$ws.Cells.Item(7,1).Value = 5
$ws.Cells.Item(7,2).Value = 88
$ws.Cells.Item(8,1).Value = 6
$ws.Cells.Item(8,2).Value = 90
$ws.Cells.Item(9,1).Value = 7
$ws.Cells.Item(9,2).Value = 93
$ws.Cells.Item(10,1).Value = 8
$ws.Cells.Item(10,2).Value = 96
$ws.Cells.Item(11,1).Value = 9
$ws.Cells.Item(11,2).Value = 99
$ws.Cells.Item(12,1).Value = 10
$ws.Cells.Item(12,2).Value = 102
$ws.Cells.Item(13,1).Value = 11
$ws.Cells.Item(13,2).Value = 105
$ws.Cells.Item(14,1).Value = 12
$ws.Cells.Item(14,2).Value = 108
$ws.Cells.Item(15,1).Value = 13
$ws.Cells.Item(15,2).Value = 111
$ws.Cells.Item(16,1).Value = 14
$ws.Cells.Item(16,2).Value = 114
$ws.Cells.Item(17,1).Value = 15
$ws.Cells.Item(17,2).Value = 117
$ws.Cells.Item(18,1).Value = 16
$ws.Cells.Item(18,2).Value = 120
$ws.Cells.Item(19,1).Value = 17
$ws.Cells.Item(19,2).Value = 123
$ws.Cells.Item(20,1).Value = 18
$ws.Cells.Item(20,2).Value = 126
$ws.Cells.Item(21,1).Value = 19
$ws.Cells.Item(21,2).Value = 129
$ws.Cells.Item(22,1).Value = 20
$ws.Cells.Item(22,2).Value = 130
$ws.Cells.Item(23,1).Value = 21
$ws.Cells.Item(23,2).Value = 130
$ws.Cells.Item(24,1).Value = 22
$ws.Cells.Item(24,2).Value = 130
$ws.Cells.Item(25,1).Value = 23
$ws.Cells.Item(25,2).Value = 130
$ws.Cells.Item(26,1).Value = 24
$ws.Cells.Item(26,2).Value = 130
$ws.Cells.Item(27,1).Value = 25
$ws.Cells.Item(27,2).Value = 130
$ws.Cells.Item(28,1).Value = 26
$ws.Cells.Item(28,2).Value = 130
$ws.Cells.Item(29,1).Value = 27
$ws.Cells.Item(29,2).Value = 130
$ws.Cells.Item(30,1).Value = 28
$ws.Cells.Item(30,2).Value = 130
$ws.Cells.Item(31,1).Value = 29
$ws.Cells.Item(31,2).Value = 130
$ws.Cells.Item(32,1).Value = 30
$ws.Cells.Item(32,2).Value = 130
$ws.Cells.Item(33,1).Value = 31
$ws.Cells.Item(33,2).Value = 130
$ws.Cells.Item(34,1).Value = 32
$ws.Cells.Item(34,2).Value = 130

$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(2,1).Value = 0
$ws.Cells.Item(2,2).Value = 70
$ws.Cells.Item(3,1).Value = 1
$ws.Cells.Item(3,2).Value = 72
$ws.Cells.Item(4,1).Value = 2
$ws.Cells.Item(4,2).Value = 75
$ws.Cells.Item(5,1).Value = 3
$ws.Cells.Item(5,2).Value = 78
$ws.Cells.Item(6,1).Value = 4
$ws.Cells.Item(6,2).Value = 81
$ws.Cells.Item(7,1).Value = 5
$ws.Cells.Item(7,2).Value = 83
$ws.Cells.Item(8,1).Value = 6
$ws.Cells.Item(8,2).Value = 86
$ws.Cells.Item(9,1).Value = 7
$ws.Cells.Item(9,2).Value = 89
$ws.Cells.Item(10,1).Value = 8
$ws.Cells.Item(10,2).Value = 92
$ws.Cells.Item(11,1).Value = 9
$ws.Cells.Item(11,2).Value = 95
$ws.Cells.Item(12,1).Value = 10
$ws.Cells.Item(12,2).Value = 97
$ws.Cells.Item(13,1).Value = 11
$ws.Cells.Item(13,2).Value = 100
$ws.Cells.Item(14,1).Value = 12
$ws.Cells.Item(14,2).Value = 103
$ws.Cells.Item(15,1).Value = 13
$ws.Cells.Item(15,2).Value = 106
$ws.Cells.Item(16,1).Value = 14
$ws.Cells.Item(16,2).Value = 109
$ws.Cells.Item(17,1).Value = 15
$ws.Cells.Item(17,2).Value = 112
$ws.Cells.Item(18,1).Value = 16
$ws.Cells.Item(18,2).Value = 115
$ws.Cells.Item(19,1).Value = 17
$ws.Cells.Item(19,2).Value = 118
$ws.Cells.Item(20,1).Value = 18
$ws.Cells.Item(20,2).Value = 121
$ws.Cells.Item(21,1).Value = 19
$ws.Cells.Item(21,2).Value = 124
$ws.Cells.Item(22,1).Value = 20
$ws.Cells.Item(22,2).Value = 127
$ws.Cells.Item(23,1).Value = 21
$ws.Cells.Item(23,2).Value = 130
$ws.Cells.Item(24,1).Value = 22
$ws.Cells.Item(24,2).Value = 130
$ws.Cells.Item(25,1).Value = 23
$ws.Cells.Item(25,2).Value = 130
$ws.Cells.Item(26,1).Value = 24
$ws.Cells.Item(26,2).Value = 130
$ws.Cells.Item(27,1).Value = 25
$ws.Cells.Item(27,2).Value = 130
$ws.Cells.Item(28,1).Value = 26
$ws.Cells.Item(28,2).Value = 130
$ws.Cells.Item(29,1).Value = 27
$ws.Cells.Item(29,2).Value = 130
$ws.Cells.Item(30,1).Value = 28
$ws.Cells.Item(30,2).Value = 130
$ws.Cells.Item(31,1).Value = 29
$ws.Cells.Item(31,2).Value = 130
$ws.Cells.Item(32,1).Value = 30
$ws.Cells.Item(32,2).Value = 130
$ws.Cells.Item(33,1).Value = 31
$ws.Cells.Item(33,2).Value = 130
$ws.Cells.Item(34,1).Value = 32
$ws.Cells.Item(34,2).Value = 130

$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(2,1).Value = 0
$ws.Cells.Item(2,2).Value = 66
$ws.Cells.Item(3,1).Value = 1
$ws.Cells.Item(3,2).Value = 69
$ws.Cells.Item(4,1).Value = 2
$ws.Cells.Item(4,2).Value = 71
$ws.Cells.Item(5,1).Value = 3
$ws.Cells.Item(5,2).Value = 74
$ws.Cells.Item(6,1).Value = 4
$ws.Cells.Item(6,2).Value = 77
$ws.Cells.Item(7,1).Value = 5
$ws.Cells.Item(7,2).Value = 79
$ws.Cells.Item(8,1).Value = 6
$ws.Cells.Item(8,2).Value = 82
$ws.Cells.Item(9,1).Value = 7
$ws.Cells.Item(9,2).Value = 85
$ws.Cells.Item(10,1).Value = 8
$ws.Cells.Item(10,2).Value = 87
$ws.Cells.Item(11,1).Value = 9
$ws.Cells.Item(11,2).Value = 90
$ws.Cells.Item(12,1).Value = 10
$ws.Cells.Item(12,2).Value = 93
$ws.Cells.Item(13,1).Value = 11
$ws.Cells.Item(13,2).Value = 96
$ws.Cells.Item(14,1).Value = 12
$ws.Cells.Item(14,2).Value = 98
$ws.Cells.Item(15,1).Value = 13
$ws.Cells.Item(15,2).Value = 101
$ws.Cells.Item(16,1).Value = 14
$ws.Cells.Item(16,2).Value = 104
$ws.Cells.Item(17,1).Value = 15
$ws.Cells.Item(17,2).Value = 107
$ws.Cells.Item(18,1).Value = 16
$ws.Cells.Item(18,2).Value = 110
$ws.Cells.Item(19,1).Value = 17
$ws.Cells.Item(19,2).Value = 113
$ws.Cells.Item(20,1).Value = 18
$ws.Cells.Item(20,2).Value = 115
$ws.Cells.Item(21,1).Value = 19
$ws.Cells.Item(21,2).Value = 118
$ws.Cells.Item(22,1).Value = 20
$ws.Cells.Item(22,2).Value = 121
$ws.Cells.Item(23,1).Value = 21
$ws.Cells.Item(23,2).Value = 124
$ws.Cells.Item(24,1).Value = 22
$ws.Cells.Item(24,2).Value = 127
$ws.Cells.Item(25,1).Value = 23
$ws.Cells.Item(25,2).Value = 130
$ws.Cells.Item(26,1).Value = 24
$ws.Cells.Item(26,2).Value = 130
$ws.Cells.Item(27,1).Value = 25
$ws.Cells.Item(27,2).Value = 130
$ws.Cells.Item(28,1).Value = 26
$ws.Cells.Item(28,2).Value = 130
$ws.Cells.Item(29,1).Value = 27
$ws.Cells.Item(29,2).Value = 130
$ws.Cells.Item(30,1).Value = 28
$ws.Cells.Item(30,2).Value = 130
$ws.Cells.Item(31,1).Value = 29
$ws.Cells.Item(31,2).Value = 130
$ws.Cells.Item(32,1).Value = 30
$ws.Cells.Item(32,2).Value = 130
$ws.Cells.Item(33,1).Value = 31
$ws.Cells.Item(33,2).Value = 130
$ws.Cells.Item(34,1).Value = 32
$ws.Cells.Item(34,2).Value = 130

$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(2,1).Value = 0
$ws.Cells.Item(2,2).Value = 62
$ws.Cells.Item(3,1).Value = 1
$ws.Cells.Item(3,2).Value = 65
$ws.Cells.Item(4,1).Value = 2
$ws.Cells.Item(4,2).Value = 67
$ws.Cells.Item(5,1).Value = 3
$ws.Cells.Item(5,2).Value = 70
$ws.Cells.Item(6,1).Value = 4
$ws.Cells.Item(6,2).Value = 72
$ws.Cells.Item(7,1).Value = 5
$ws.Cells.Item(7,2).Value = 75
$ws.Cells.Item(8,1).Value = 6
$ws.Cells.Item(8,2).Value = 77
$ws.Cells.Item(9,1).Value = 7
$ws.Cells.Item(9,2).Value = 80
$ws.Cells.Item(10,1).Value = 8
$ws.Cells.Item(10,2).Value = 82
$ws.Cells.Item(11,1).Value = 9
$ws.Cells.Item(11,2).Value = 85
$ws.Cells.Item(12,1).Value = 10
$ws.Cells.Item(12,2).Value = 88
$ws.Cells.Item(13,1).Value = 11
$ws.Cells.Item(13,2).Value = 90
$ws.Cells.Item(14,1).Value = 12
$ws.Cells.Item(14,2).Value = 93
$ws.Cells.Item(15,1).Value = 13
$ws.Cells.Item(15,2).Value = 96
$ws.Cells.Item(16,1).Value = 14
$ws.Cells.Item(16,2).Value = 98
$ws.Cells.Item(17,1).Value = 15
$ws.Cells.Item(17,2).Value = 101
$ws.Cells.Item(18,1).Value = 16
$ws.Cells.Item(18,2).Value = 104
$ws.Cells.Item(19,1).Value = 17
$ws.Cells.Item(19,2).Value = 107
$ws.Cells.Item(20,1).Value = 18
$ws.Cells.Item(20,2).Value = 109
$ws.Cells.Item(21,1).Value = 19
$ws.Cells.Item(21,2).Value = 112
$ws.Cells.Item(22,1).Value = 20
$ws.Cells.Item(22,2).Value = 115
$ws.Cells.Item(23,1).Value = 21
$ws.Cells.Item(23,2).Value = 118
$ws.Cells.Item(24,1).Value = 22
$ws.Cells.Item(24,2).Value = 121
$ws.Cells.Item(25,1).Value = 23
$ws.Cells.Item(25,2).Value = 124
$ws.Cells.Item(26,1).Value = 24
$ws.Cells.Item(26,2).Value = 127
$ws.Cells.Item(27,1).Value = 25
$ws.Cells.Item(27,2).Value = 129
$ws.Cells.Item(28,1).Value = 26
$ws.Cells.Item(28,2).Value = 130
$ws.Cells.Item(29,1).Value = 27
$ws.Cells.Item(29,2).Value = 130
$ws.Cells.Item(30,1).Value = 28
$ws.Cells.Item(30,2).Value = 130
$ws.Cells.Item(31,1).Value = 29
$ws.Cells.Item(31,2).Value = 130
$ws.Cells.Item(32,1).Value = 30
$ws.Cells.Item(32,2).Value = 130
$ws.Cells.Item(33,1).Value = 31
$ws.Cells.Item(33,2).Value = 130
$ws.Cells.Item(34,1).Value = 32
$ws.Cells.Item(34,2).Value = 130

$ws = $wb.Worksheets.Item(5)
$ws.Cells.Item(2,1).Value = 0
$ws.Cells.Item(2,2).Value = 58
$ws.Cells.Item(3,1).Value = 1
$ws.Cells.Item(3,2).Value = 60
$ws.Cells.Item(4,1).Value = 2
$ws.Cells.Item(4,2).Value = 62
$ws.Cells.Item(5,1).Value = 3
$ws.Cells.Item(5,2).Value = 65
$ws.Cells.Item(6,1).Value = 4
$ws.Cells.Item(6,2).Value = 67
$ws.Cells.Item(7,1).Value = 5
$ws.Cells.Item(7,2).Value = 69
$ws.Cells.Item(8,1).Value = 6
$ws.Cells.Item(8,2).Value = 72
$ws.Cells.Item(9,1).Value = 7
$ws.Cells.Item(9,2).Value = 74
$ws.Cells.Item(10,1).Value = 8
$ws.Cells.Item(10,2).Value = 77
$ws.Cells.Item(11,1).Value = 9
$ws.Cells.Item(11,2).Value = 79
$ws.Cells.Item(12,1).Value = 10
$ws.Cells.Item(12,2).Value = 82
$ws.Cells.Item(13,1).Value = 11
$ws.Cells.Item(13,2).Value = 84
$ws.Cells.Item(14,1).Value = 12
$ws.Cells.Item(14,2).Value = 87
$ws.Cells.Item(15,1).Value = 13
$ws.Cells.Item(15,2).Value = 89
$ws.Cells.Item(16,1).Value = 14
$ws.Cells.Item(16,2).Value = 92
$ws.Cells.Item(17,1).Value = 15
$ws.Cells.Item(17,2).Value = 95
$ws.Cells.Item(18,1).Value = 16
$ws.Cells.Item(18,2).Value = 97
$ws.Cells.Item(19,1).Value = 17
$ws.Cells.Item(19,2).Value = 100
$ws.Cells.Item(20,1).Value = 18
$ws.Cells.Item(20,2).Value = 103
$ws.Cells.Item(21,1).Value = 19
$ws.Cells.Item(21,2).Value = 105
$ws.Cells.Item(22,1).Value = 20
$ws.Cells.Item(22,2).Value = 108
$ws.Cells.Item(23,1).Value = 21
$ws.Cells.Item(23,2).Value = 111
$ws.Cells.Item(24,1).Value = 22
$ws.Cells.Item(24,2).Value = 114
$ws.Cells.Item(25,1).Value = 23
$ws.Cells.Item(25,2).Value = 117
$ws.Cells.Item(26,1).Value = 24
$ws.Cells.Item(26,2).Value = 119
$ws.Cells.Item(27,1).Value = 25
$ws.Cells.Item(27,2).Value = 122
$ws.Cells.Item(28,1).Value = 26
$ws.Cells.Item(28,2).Value = 125
$ws.Cells.Item(29,1).Value = 27
$ws.Cells.Item(29,2).Value = 128
$ws.Cells.Item(30,1).Value = 28
$ws.Cells.Item(30,2).Value = 130
$ws.Cells.Item(31,1).Value = 29
$ws.Cells.Item(31,2).Value = 130
$ws.Cells.Item(32,1).Value = 30
$ws.Cells.Item(32,2).Value = 130
$ws.Cells.Item(33,1).Value = 31
$ws.Cells.Item(33,2).Value = 130
$ws.Cells.Item(34,1).Value = 32
$ws.Cells.Item(34,2).Value = 130

$ws = $wb.Worksheets.Item(6)
$ws.Cells.Item(2,1).Value = 0
$ws.Cells.Item(2,2).Value = 53
$ws.Cells.Item(3,1).Value = 1
$ws.Cells.Item(3,2).Value = 56
$ws.Cells.Item(4,1).Value = 2
$ws.Cells.Item(4,2).Value = 58
$ws.Cells.Item(5,1).Value = 3
$ws.Cells.Item(5,2).Value = 60
$ws.Cells.Item(6,1).Value = 4
$ws.Cells.Item(6,2).Value = 62
$ws.Cells.Item(7,1).Value = 5
$ws.Cells.Item(7,2).Value = 65
$ws.Cells.Item(8,1).Value = 6
$ws.Cells.Item(8,2).Value = 67
$ws.Cells.Item(9,1).Value = 7
$ws.Cells.Item(9,2).Value = 69
$ws.Cells.Item(10,1).Value = 8
$ws.Cells.Item(10,2).Value = 72
$ws.Cells.Item(11,1).Value = 9
$ws.Cells.Item(11,2).Value = 74
$ws.Cells.Item(12,1).Value = 10
$ws.Cells.Item(12,2).Value = 76
$ws.Cells.Item(13,1).Value = 11
$ws.Cells.Item(13,2).Value = 79
$ws.Cells.Item(14,1).Value = 12
$ws.Cells.Item(14,2).Value = 81
$ws.Cells.Item(15,1).Value = 13
$ws.Cells.Item(15,2).Value = 84
$ws.Cells.Item(16,1).Value = 14
$ws.Cells.Item(16,2).Value = 86
$ws.Cells.Item(17,1).Value = 15
$ws.Cells.Item(17,2).Value = 89
$ws.Cells.Item(18,1).Value = 16
$ws.Cells.Item(18,2).Value = 91
$ws.Cells.Item(19,1).Value = 17
$ws.Cells.Item(19,2).Value = 94
$ws.Cells.Item(20,1).Value = 18
$ws.Cells.Item(20,2).Value = 97
$ws.Cells.Item(21,1).Value = 19
$ws.Cells.Item(21,2).Value = 99
$ws.Cells.Item(22,1).Value = 20
$ws.Cells.Item(22,2).Value = 102
$ws.Cells.Item(23,1).Value = 21
$ws.Cells.Item(23,2).Value = 105
$ws.Cells.Item(24,1).Value = 22
$ws.Cells.Item(24,2).Value = 107
$ws.Cells.Item(25,1).Value = 23
$ws.Cells.Item(25,2).Value = 110
$ws.Cells.Item(26,1).Value = 24
$ws.Cells.Item(26,2).Value = 113
$ws.Cells.Item(27,1).Value = 25
$ws.Cells.Item(27,2).Value = 116
$ws.Cells.Item(28,1).Value = 26
$ws.Cells.Item(28,2).Value = 119
$ws.Cells.Item(29,1).Value = 27
$ws.Cells.Item(29,2).Value = 122
$ws.Cells.Item(30,1).Value = 28
$ws.Cells.Item(30,2).Value = 125
$ws.Cells.Item(31,1).Value = 29
$ws.Cells.Item(31,2).Value = 128
$ws.Cells.Item(32,1).Value = 30
$ws.Cells.Item(32,2).Value = 130
$ws.Cells.Item(33,1).Value = 31
$ws.Cells.Item(33,2).Value = 130
$ws.Cells.Item(34,1).Value = 32
$ws.Cells.Item(34,2).Value = 130

$ws = $wb.Worksheets.Item(7)
$ws.Cells.Item(2,1).Value = 0
$ws.Cells.Item(2,2).Value = 50
$ws.Cells.Item(3,1).Value = 1
$ws.Cells.Item(3,2).Value = 52
$ws.Cells.Item(4,1).Value = 2
$ws.Cells.Item(4,2).Value = 54
$ws.Cells.Item(5,1).Value = 3
$ws.Cells.Item(5,2).Value = 56
$ws.Cells.Item(6,1).Value = 4
$ws.Cells.Item(6,2).Value = 58
$ws.Cells.Item(7,1).Value = 5
$ws.Cells.Item(7,2).Value = 60
$ws.Cells.Item(8,1).Value = 6
$ws.Cells.Item(8,2).Value = 62
$ws.Cells.Item(9,1).Value = 7
$ws.Cells.Item(9,2).Value = 65
$ws.Cells.Item(10,1).Value = 8
$ws.Cells.Item(10,2).Value = 67
$ws.Cells.Item(11,1).Value = 9
$ws.Cells.Item(11,2).Value = 69
$ws.Cells.Item(12,1).Value = 10
$ws.Cells.Item(12,2).Value = 71
$ws.Cells.Item(13,1).Value = 11
$ws.Cells.Item(13,2).Value = 74
$ws.Cells.Item(14,1).Value = 12
$ws.Cells.Item(14,2).Value = 76
$ws.Cells.Item(15,1).Value = 13
$ws.Cells.Item(15,2).Value = 79
$ws.Cells.Item(16,1).Value = 14
$ws.Cells.Item(16,2).Value = 81
$ws.Cells.Item(17,1).Value = 15
$ws.Cells.Item(17,2).Value = 83
$ws.Cells.Item(18,1).Value = 16
$ws.Cells.Item(18,2).Value = 86
$ws.Cells.Item(19,1).Value = 17
$ws.Cells.Item(19,2).Value = 88
$ws.Cells.Item(20,1).Value = 18
$ws.Cells.Item(20,2).Value = 91
$ws.Cells.Item(21,1).Value = 19
$ws.Cells.Item(21,2).Value = 94
$ws.Cells.Item(22,1).Value = 20
$ws.Cells.Item(22,2).Value = 96
$ws.Cells.Item(23,1).Value = 21
$ws.Cells.Item(23,2).Value = 99
$ws.Cells.Item(24,1).Value = 22
$ws.Cells.Item(24,2).Value = 102
$ws.Cells.Item(25,1).Value = 23
$ws.Cells.Item(25,2).Value = 104
$ws.Cells.Item(26,1).Value = 24
$ws.Cells.Item(26,2).Value = 107
$ws.Cells.Item(27,1).Value = 25
$ws.Cells.Item(27,2).Value = 110
$ws.Cells.Item(28,1).Value = 26
$ws.Cells.Item(28,2).Value = 113
$ws.Cells.Item(29,1).Value = 27
$ws.Cells.Item(29,2).Value = 116
$ws.Cells.Item(30,1).Value = 28
$ws.Cells.Item(30,2).Value = 119
$ws.Cells.Item(31,1).Value = 29
$ws.Cells.Item(31,2).Value = 123
$ws.Cells.Item(32,1).Value = 30
$ws.Cells.Item(32,2).Value = 126
$ws.Cells.Item(33,1).Value = 31
$ws.Cells.Item(33,2).Value = 129
$ws.Cells.Item(34,1).Value = 32
$ws.Cells.Item(34,2).Value = 130

$ws = $wb.Worksheets.Item(8)
$ws.Cells.Item(2,1).Value = 0
$ws.Cells.Item(2,2).Value = 46
$ws.Cells.Item(3,1).Value = 1
$ws.Cells.Item(3,2).Value = 48
$ws.Cells.Item(4,1).Value = 2
$ws.Cells.Item(4,2).Value = 50
$ws.Cells.Item(5,1).Value = 3
$ws.Cells.Item(5,2).Value = 52
$ws.Cells.Item(6,1).Value = 4
$ws.Cells.Item(6,2).Value = 54
$ws.Cells.Item(7,1).Value = 5
$ws.Cells.Item(7,2).Value = 56
$ws.Cells.Item(8,1).Value = 6
$ws.Cells.Item(8,2).Value = 58
$ws.Cells.Item(9,1).Value = 7
$ws.Cells.Item(9,2).Value = 61
$ws.Cells.Item(10,1).Value = 8
$ws.Cells.Item(10,2).Value = 63
$ws.Cells.Item(11,1).Value = 9
$ws.Cells.Item(11,2).Value = 65
$ws.Cells.Item(12,1).Value = 10
$ws.Cells.Item(12,2).Value = 67
$ws.Cells.Item(13,1).Value = 11
$ws.Cells.Item(13,2).Value = 69
$ws.Cells.Item(14,1).Value = 12
$ws.Cells.Item(14,2).Value = 72
$ws.Cells.Item(15,1).Value = 13
$ws.Cells.Item(15,2).Value = 74
$ws.Cells.Item(16,1).Value = 14
$ws.Cells.Item(16,2).Value = 76
$ws.Cells.Item(17,1).Value = 15
$ws.Cells.Item(17,2).Value = 79
$ws.Cells.Item(18,1).Value = 16
$ws.Cells.Item(18,2).Value = 81
$ws.Cells.Item(19,1).Value = 17
$ws.Cells.Item(19,2).Value = 83
$ws.Cells.Item(20,1).Value = 18
$ws.Cells.Item(20,2).Value = 86
$ws.Cells.Item(21,1).Value = 19
$ws.Cells.Item(21,2).Value = 89
$ws.Cells.Item(22,1).Value = 20
$ws.Cells.Item(22,2).Value = 91
$ws.Cells.Item(23,1).Value = 21
$ws.Cells.Item(23,2).Value = 94
$ws.Cells.Item(24,1).Value = 22
$ws.Cells.Item(24,2).Value = 97
$ws.Cells.Item(25,1).Value = 23
$ws.Cells.Item(25,2).Value = 99
$ws.Cells.Item(26,1).Value = 24
$ws.Cells.Item(26,2).Value = 102
$ws.Cells.Item(27,1).Value = 25
$ws.Cells.Item(27,2).Value = 105
$ws.Cells.Item(28,1).Value = 26
$ws.Cells.Item(28,2).Value = 108
$ws.Cells.Item(29,1).Value = 27
$ws.Cells.Item(29,2).Value = 111
$ws.Cells.Item(30,1).Value = 28
$ws.Cells.Item(30,2).Value = 115
$ws.Cells.Item(31,1).Value = 29
$ws.Cells.Item(31,2).Value = 118
$ws.Cells.Item(32,1).Value = 30
$ws.Cells.Item(32,2).Value = 122
$ws.Cells.Item(33,1).Value = 31
$ws.Cells.Item(33,2).Value = 125
$ws.Cells.Item(34,1).Value = 32
$ws.Cells.Item(34,2).Value = 129

$ws = $wb.Worksheets.Item(9)
$ws.Cells.Item(2,1).Value = 0
$ws.Cells.Item(2,2).Value = 42
$ws.Cells.Item(3,1).Value = 1
$ws.Cells.Item(3,2).Value = 44
$ws.Cells.Item(4,1).Value = 2
$ws.Cells.Item(4,2).Value = 46
$ws.Cells.Item(5,1).Value = 3
$ws.Cells.Item(5,2).Value = 48
$ws.Cells.Item(6,1).Value = 4
$ws.Cells.Item(6,2).Value = 50
$ws.Cells.Item(7,1).Value = 5
$ws.Cells.Item(7,2).Value = 52
$ws.Cells.Item(8,1).Value = 6
$ws.Cells.Item(8,2).Value = 53
$ws.Cells.Item(9,1).Value = 7
$ws.Cells.Item(9,2).Value = 55
$ws.Cells.Item(10,1).Value = 8
$ws.Cells.Item(10,2).Value = 58
$ws.Cells.Item(11,1).Value = 9
$ws.Cells.Item(11,2).Value = 60
$ws.Cells.Item(12,1).Value = 10
$ws.Cells.Item(12,2).Value = 62
$ws.Cells.Item(13,1).Value = 11
$ws.Cells.Item(13,2).Value = 64
$ws.Cells.Item(14,1).Value = 12
$ws.Cells.Item(14,2).Value = 66
$ws.Cells.Item(15,1).Value = 13
$ws.Cells.Item(15,2).Value = 68
$ws.Cells.Item(16,1).Value = 14
$ws.Cells.Item(16,2).Value = 70
$ws.Cells.Item(17,1).Value = 15
$ws.Cells.Item(17,2).Value = 73
$ws.Cells.Item(18,1).Value = 16
$ws.Cells.Item(18,2).Value = 75
$ws.Cells.Item(19,1).Value = 17
$ws.Cells.Item(19,2).Value = 78
$ws.Cells.Item(20,1).Value = 18
$ws.Cells.Item(20,2).Value = 80
$ws.Cells.Item(21,1).Value = 19
$ws.Cells.Item(21,2).Value = 83
$ws.Cells.Item(22,1).Value = 20
$ws.Cells.Item(22,2).Value = 85
$ws.Cells.Item(23,1).Value = 21
$ws.Cells.Item(23,2).Value = 88
$ws.Cells.Item(24,1).Value = 22
$ws.Cells.Item(24,2).Value = 91
$ws.Cells.Item(25,1).Value = 23
$ws.Cells.Item(25,2).Value = 94
$ws.Cells.Item(26,1).Value = 24
$ws.Cells.Item(26,2).Value = 97
$ws.Cells.Item(27,1).Value = 25
$ws.Cells.Item(27,2).Value = 100
$ws.Cells.Item(28,1).Value = 26
$ws.Cells.Item(28,2).Value = 103
$ws.Cells.Item(29,1).Value = 27
$ws.Cells.Item(29,2).Value = 107
$ws.Cells.Item(30,1).Value = 28
$ws.Cells.Item(30,2).Value = 110
$ws.Cells.Item(31,1).Value = 29
$ws.Cells.Item(31,2).Value = 114
$ws.Cells.Item(32,1).Value = 30
$ws.Cells.Item(32,2).Value = 118
$ws.Cells.Item(33,1).Value = 31
$ws.Cells.Item(33,2).Value = 123
$ws.Cells.Item(34,1).Value = 32
$ws.Cells.Item(34,2).Value = 127
